# Auto-generated Excel COM-interop edit script
# Applies cell-value updates to the Raiden_Profits workbook (market-data refresh)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ALC_updates = @(
    @{ Cell = "H80"; Value = 299.6154 }
    @{ Cell = "J80"; Value = 165.90909 }
    @{ Cell = "L80"; Value = 497.72727 }
    @{ Cell = "N80"; Value = -2493.72727 }
    @{ Cell = "H83"; Value = 299.6154 }
    @{ Cell = "J83"; Value = 165.90909 }
    @{ Cell = "L83"; Value = 1493.18181 }
    @{ Cell = "N83"; Value = -11477.18181 }
    @{ Cell = "H96"; Value = 1175.6 }
    @{ Cell = "I96"; Value = 2725 }
    @{ Cell = "J96"; Value = 142.66667 }
    @{ Cell = "K96"; Value = 8175 }
    @{ Cell = "L96"; Value = 428.00001 }
    @{ Cell = "M96"; Value = -6802 }
    @{ Cell = "N96"; Value = -3174.00001 }
    @{ Cell = "H112"; Value = 1470.75 }
    @{ Cell = "J112"; Value = 1925.1666 }
    @{ Cell = "L112"; Value = 5775.4998 }
    @{ Cell = "N112"; Value = -7991.4998 }
    @{ Cell = "H116"; Value = 4913.933 }
    @{ Cell = "I116"; Value = 5063 }
    @{ Cell = "K116"; Value = 5063 }
    @{ Cell = "M116"; Value = -1621 }
    @{ Cell = "H132"; Value = 3209.5557 }
    @{ Cell = "I132"; Value = 2986.5881 }
    @{ Cell = "J132"; Value = 7000 }
    @{ Cell = "K132"; Value = 8959.764299999999 }
    @{ Cell = "L132"; Value = 21000 }
    @{ Cell = "M132"; Value = -6429.764299999999 }
    @{ Cell = "N132"; Value = -26060 }
    @{ Cell = "H138"; Value = 2471.2292 }
    @{ Cell = "J138"; Value = 2277.5527 }
    @{ Cell = "L138"; Value = 6832.658100000001 }
    @{ Cell = "N138"; Value = -17112.6581 }
    @{ Cell = "H141"; Value = 5106.6924 }
    @{ Cell = "J141"; Value = 9999.5 }
    @{ Cell = "L141"; Value = 29998.5 }
    @{ Cell = "N141"; Value = -40358.5 }
)
foreach ($u in $ALC_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ARM_updates = @(
    @{ Cell = "H45"; Value = 1848.8572 }
    @{ Cell = "I45"; Value = 1643.7778 }
    @{ Cell = "J45"; Value = 2218 }
    @{ Cell = "K45"; Value = 1643.7778 }
    @{ Cell = "L45"; Value = 2218 }
    @{ Cell = "M45"; Value = -1266.7778 }
    @{ Cell = "N45"; Value = -2972 }
    @{ Cell = "H102"; Value = 2949.4 }
    @{ Cell = "I102"; Value = 2949.4 }
    @{ Cell = "K102"; Value = 2949.4 }
    @{ Cell = "M102"; Value = -1327.4 }
    @{ Cell = "H132"; Value = 3562.5 }
    @{ Cell = "I132"; Value = 3293.1667 }
    @{ Cell = "K132"; Value = 9879.500100000001 }
    @{ Cell = "M132"; Value = -7349.500100000001 }
)
foreach ($u in $ARM_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$BSM_updates = @(
    @{ Cell = "H40"; Value = 0 }
    @{ Cell = "J40"; Value = 0 }
    @{ Cell = "L40"; Value = 0 }
    @{ Cell = "H82"; Value = 28691.834 }
    @{ Cell = "I82"; Value = 14084 }
    @{ Cell = "J82"; Value = 43299.668 }
    @{ Cell = "K82"; Value = 14084 }
    @{ Cell = "L82"; Value = 43299.668 }
    @{ Cell = "M82"; Value = -13701 }
    @{ Cell = "N82"; Value = -44065.668 }
    @{ Cell = "H85"; Value = 28691.834 }
    @{ Cell = "I85"; Value = 14084 }
    @{ Cell = "J85"; Value = 43299.668 }
    @{ Cell = "K85"; Value = 14084 }
    @{ Cell = "L85"; Value = 43299.668 }
    @{ Cell = "M85"; Value = -12758 }
    @{ Cell = "N85"; Value = -45951.668 }
    @{ Cell = "H86"; Value = 3204.4 }
    @{ Cell = "I86"; Value = 4000.6 }
    @{ Cell = "J86"; Value = 2408.2 }
    @{ Cell = "K86"; Value = 4000.6 }
    @{ Cell = "L86"; Value = 2408.2 }
    @{ Cell = "M86"; Value = -2877.6 }
    @{ Cell = "N86"; Value = -4654.2 }
    @{ Cell = "H89"; Value = 3204.4 }
    @{ Cell = "I89"; Value = 4000.6 }
    @{ Cell = "J89"; Value = 2408.2 }
    @{ Cell = "K89"; Value = 20003 }
    @{ Cell = "L89"; Value = 12041 }
    @{ Cell = "M89"; Value = -14387 }
    @{ Cell = "N89"; Value = -23273 }
    @{ Cell = "H94"; Value = 2369.842 }
    @{ Cell = "I94"; Value = 619.35297 }
    @{ Cell = "K94"; Value = 619.35297 }
    @{ Cell = "M94"; Value = -168.35297 }
    @{ Cell = "H96"; Value = 13133.556 }
    @{ Cell = "I96"; Value = 13133.556 }
    @{ Cell = "J96"; Value = 0 }
    @{ Cell = "K96"; Value = 13133.556 }
    @{ Cell = "L96"; Value = 0 }
    @{ Cell = "M96"; Value = -10387.556 }
)
foreach ($u in $BSM_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
$BSM_deletes = @("N40", "N96")
foreach ($d in $BSM_deletes) {
    $ws.Range($d).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$CRP_updates = @(
    @{ Cell = "H130"; Value = 58726.5 }
)
foreach ($u in $CRP_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$CUL_updates = @(
    @{ Cell = "H2"; Value = 1000351.6 }
    @{ Cell = "J2"; Value = 419.66666 }
    @{ Cell = "L2"; Value = 2517.99996 }
    @{ Cell = "N2"; Value = -2743.99996 }
    @{ Cell = "H23"; Value = 443.33334 }
    @{ Cell = "I23"; Value = 394.5 }
    @{ Cell = "K23"; Value = 1183.5 }
    @{ Cell = "M23"; Value = -948.5 }
    @{ Cell = "H26"; Value = 417.5 }
    @{ Cell = "I26"; Value = 476.5 }
    @{ Cell = "J26"; Value = 299.5 }
    @{ Cell = "K26"; Value = 1429.5 }
    @{ Cell = "L26"; Value = 898.5 }
    @{ Cell = "M26"; Value = -1141.5 }
    @{ Cell = "N26"; Value = -1474.5 }
    @{ Cell = "H34"; Value = 3523.5 }
    @{ Cell = "I34"; Value = 195 }
    @{ Cell = "J34"; Value = 4633 }
    @{ Cell = "K34"; Value = 585 }
    @{ Cell = "L34"; Value = 13899 }
    @{ Cell = "M34"; Value = -501 }
    @{ Cell = "N34"; Value = -14067 }
    @{ Cell = "H39"; Value = 5904.7334 }
    @{ Cell = "J39"; Value = 6627.385 }
    @{ Cell = "L39"; Value = 19882.155 }
    @{ Cell = "N39"; Value = -20470.155 }
    @{ Cell = "H81"; Value = 3998.5 }
    @{ Cell = "I81"; Value = 3998.5 }
    @{ Cell = "K81"; Value = 11995.5 }
    @{ Cell = "M81"; Value = -10872.5 }
    @{ Cell = "H84"; Value = 3998.5 }
    @{ Cell = "I84"; Value = 3998.5 }
    @{ Cell = "K84"; Value = 35986.5 }
    @{ Cell = "M84"; Value = -30370.5 }
    @{ Cell = "H108"; Value = 2145 }
    @{ Cell = "I108"; Value = 2145 }
    @{ Cell = "K108"; Value = 6435 }
    @{ Cell = "M108"; Value = -3555 }
    @{ Cell = "H109"; Value = 235 }
    @{ Cell = "J109"; Value = 0 }
    @{ Cell = "L109"; Value = 0 }
    @{ Cell = "H114"; Value = 3310.7144 }
    @{ Cell = "I114"; Value = 2797 }
    @{ Cell = "J114"; Value = 3995.6667 }
    @{ Cell = "K114"; Value = 8391 }
    @{ Cell = "L114"; Value = 11987.0001 }
    @{ Cell = "M114"; Value = -5137 }
    @{ Cell = "N114"; Value = -18495.0001 }
    @{ Cell = "H117"; Value = 1368.5 }
    @{ Cell = "J117"; Value = 1612.1428 }
    @{ Cell = "L117"; Value = 4836.428400000001 }
    @{ Cell = "N117"; Value = -11720.4284 }
    @{ Cell = "H120"; Value = 10400 }
    @{ Cell = "I120"; Value = 4000 }
    @{ Cell = "K120"; Value = 12000 }
    @{ Cell = "M120"; Value = -7162 }
    @{ Cell = "H129"; Value = 3514.4443 }
    @{ Cell = "J129"; Value = 3642.9412 }
    @{ Cell = "L129"; Value = 10928.8236 }
    @{ Cell = "N129"; Value = -20928.8236 }
    @{ Cell = "H131"; Value = 24735.229 }
    @{ Cell = "I131"; Value = 223253.8 }
    @{ Cell = "J131"; Value = 1651.6744 }
    @{ Cell = "K131"; Value = 669761.3999999999 }
    @{ Cell = "L131"; Value = 4955.023200000001 }
    @{ Cell = "M131"; Value = -664721.3999999999 }
    @{ Cell = "N131"; Value = -15035.0232 }
    @{ Cell = "H138"; Value = 1578.8 }
    @{ Cell = "I138"; Value = 1578.8 }
    @{ Cell = "J138"; Value = 0 }
    @{ Cell = "K138"; Value = 4736.4 }
    @{ Cell = "L138"; Value = 0 }
    @{ Cell = "M138"; Value = 403.6000000000004 }
    @{ Cell = "H140"; Value = 3431.55 }
    @{ Cell = "I140"; Value = 2142.611 }
    @{ Cell = "K140"; Value = 6427.833 }
    @{ Cell = "M140"; Value = -1247.833 }
)
foreach ($u in $CUL_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
$CUL_deletes = @("N109", "N138")
foreach ($d in $CUL_deletes) {
    $ws.Range($d).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$GSM_updates = @(
    @{ Cell = "H80"; Value = 8325 }
    @{ Cell = "J80"; Value = 14422.667 }
    @{ Cell = "L80"; Value = 14422.667 }
    @{ Cell = "N80"; Value = -16418.667 }
    @{ Cell = "H83"; Value = 8325 }
    @{ Cell = "J83"; Value = 14422.667 }
    @{ Cell = "L83"; Value = 72113.33499999999 }
    @{ Cell = "N83"; Value = -82097.33499999999 }
    @{ Cell = "H122"; Value = 1414.3077 }
    @{ Cell = "I122"; Value = 1221.5454 }
    @{ Cell = "K122"; Value = 3664.6362 }
    @{ Cell = "M122"; Value = -1214.6362 }
    @{ Cell = "H135"; Value = 67500 }
    @{ Cell = "J135"; Value = 67500 }
    @{ Cell = "L135"; Value = 67500 }
    @{ Cell = "N135"; Value = -77640 }
)
foreach ($u in $GSM_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$LTW_updates = @(
    @{ Cell = "H22"; Value = 1177.6 }
    @{ Cell = "I22"; Value = 569.4 }
    @{ Cell = "J22"; Value = 1785.8 }
    @{ Cell = "K22"; Value = 569.4 }
    @{ Cell = "L22"; Value = 1785.8 }
    @{ Cell = "M22"; Value = -274.4 }
    @{ Cell = "N22"; Value = -2375.8 }
    @{ Cell = "H27"; Value = 1177.6 }
    @{ Cell = "I27"; Value = 569.4 }
    @{ Cell = "J27"; Value = 1785.8 }
    @{ Cell = "K27"; Value = 569.4 }
    @{ Cell = "L27"; Value = 1785.8 }
    @{ Cell = "M27"; Value = -462.4 }
    @{ Cell = "N27"; Value = -1999.8 }
    @{ Cell = "H46"; Value = 2521.75 }
    @{ Cell = "I46"; Value = 933.5714 }
    @{ Cell = "J46"; Value = 3376.923 }
    @{ Cell = "K46"; Value = 933.5714 }
    @{ Cell = "L46"; Value = 3376.923 }
    @{ Cell = "M46"; Value = -745.5714 }
    @{ Cell = "N46"; Value = -3752.923 }
    @{ Cell = "H55"; Value = 534.86664 }
    @{ Cell = "I55"; Value = 398.72726 }
    @{ Cell = "K55"; Value = 398.72726 }
    @{ Cell = "M55"; Value = -225.72726 }
)
foreach ($u in $LTW_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$WVR_updates = @(
    @{ Cell = "H80"; Value = 25000 }
    @{ Cell = "J80"; Value = 25000 }
    @{ Cell = "L80"; Value = 25000 }
    @{ Cell = "N80"; Value = -26996 }
    @{ Cell = "H83"; Value = 25000 }
    @{ Cell = "J83"; Value = 25000 }
    @{ Cell = "L83"; Value = 75000 }
    @{ Cell = "N83"; Value = -84984 }
    @{ Cell = "H95"; Value = 39000 }
    @{ Cell = "J95"; Value = 39000 }
    @{ Cell = "L95"; Value = 39000 }
    @{ Cell = "N95"; Value = -44492 }
    @{ Cell = "H107"; Value = 806.375 }
    @{ Cell = "I107"; Value = 864.1818 }
    @{ Cell = "K107"; Value = 2592.5454 }
    @{ Cell = "M107"; Value = -672.5454 }
    @{ Cell = "H122"; Value = 2516.8838 }
    @{ Cell = "I122"; Value = 2035.6 }
    @{ Cell = "K122"; Value = 6106.799999999999 }
    @{ Cell = "M122"; Value = -3656.799999999999 }
    @{ Cell = "H126"; Value = 3676.9412 }
    @{ Cell = "I126"; Value = 3706 }
    @{ Cell = "K126"; Value = 11118 }
    @{ Cell = "M126"; Value = -8648 }
)
foreach ($u in $WVR_updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Applied all Raiden_Profits cell updates."